$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$range = $ws.Range("A2:A27")
$range.NumberFormat = "@"
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 1).Value = "08.19.19"
}

$ws.Range("A3:A27").Select()
